$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Arreglo errores interfaz gráfica
$ws.Range("B2").Value = 8.0
$ws.Range("B4").Value = 29.0
$ws.Range("B9").Value = 0.0
$ws.Range("B10").Value = 3.0
